$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Yes/No flag values: rows 2 and 3 become "Yes", row 4 becomes "No"
$ws.Range("A2").Value = "Yes"
$ws.Range("A3").Value = "Yes"
$ws.Range("A4").Value = "No"

# Move the active selection from A4 to A2
$ws.Range("A2").Select()
